$d = $word.ActiveDocument

# --- Paragraph 1 (title): "DD.MM.YY" line + title line, separated by <w:br/> ---
$d.Paragraphs(1).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>המאמר היומי של מייק - 09.03.25</w:t><w:br/><w:t>THE SUPER WEIGHT IN LARGE LANGUAGE MODELS</w:t></w:r></w:p>")

# --- Paragraphs 2-7 (body text) ---
$d.Paragraphs(2).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>זה די לא ייאמן, אבל מודלים שפה גדולים עם מיליארדי או אפילו עשרות מיליארדי פרמטרים עלולים לסבול ירידה כואבת בביצועים אם מורידים מהם אפילו משקל בודד. ממצא מפתיע זה חל לפחות על חלק מהמודלים העוצמתיים האלה.</w:t></w:r></w:p>")
$d.Paragraphs(3).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>מאמר זה מתעמק במאפיין ספציפי ובלתי צפוי של מודלים שפה גדולים: קיומם של `"משקלים על (SWs)`". המחברים מתקדמים מעבר לתצפית ידועה על כך ש-LLMs  מכילים משקלים חריגים המשפיעים באופן ניכר על הביצועים, ומציגים ראיות לכך שמשקל בודד יכול להיות קריטי באופן לא פרופורציונלי לתפקוד המודל.</w:t></w:r></w:p>")
$d.Paragraphs(4).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>כאמור הממצא המרכזי הוא שהורדת SW בודד יכול לגרום לירידה קשה בביצועי LLM. השפעה דרסטית זו מתבטאת כעלייה חדה בפרפלקסיטי וירידה בדיוק zero-shot לרמות כמעט אקראיות. מה שראוי לציון במיוחד הוא העובדה שהסרת SW לבין ההשפעה הקטנה יחסית של הורדה של משקלים חריגים אחרים, אפילו בעלי גודל גדול יותר.</w:t></w:r></w:p>")
$d.Paragraphs(5).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>המאמר נותן דוגמה מעניינת להשפעה של הסרת משקל על-כבד עבור הפרומפט: `"קיץ חם. חורף הוא…`"(באנגלית). הטוקן הבא הנכון צריך להיות `"קר`" ועם המודל המקורי עם SW, הוא חוזה נכון את הטוקן הבא `"קר`" בהסתברות גבוהה של 81.4%. כאשר SW מוסר, החיזוי המוביל של המודל הוא stopword `"ה`"(the) בהסתברות נמוכה ולא בטוחה של 9.0%. זה מצביע על כך ש-SW חיוני למודל כדי לבצע חיזוי נכון ובטוח של מילים משמעותיות. המאמר לא רק מתעד את התופעה הזו; הוא גם בוחן את המנגנונים הבסיסיים הקשורים אליה. המחברים מקשרים SW ליצירת `"אקטיבציות SW`", שהן אקטיבציה גדולות וחריגות המתפשטות דרך המודל כמעט ללא קשר לקלט.</w:t></w:r></w:p>")
$d.Paragraphs(6).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>יתר על כן, המחקר בוחן את ההשלכות של SW עבור קוונטיזציה של מודלי שפה. נוכחותם של חריגים, כולל SW ואקטיבציות חריגות הנגרמות מהם, מציבה אתגר משמעותי לקוונטיזציה יעילה, שכן חריגים אלה יכולים לעוות את תהליך הקוונטיזציה ולהוביל לאובדן מידע משמעותי. המחברים מדגימים ששימור חריגי SW (גם משקלים וגם אקטיבציות) יכול לשפר את יעילות הקוונטיזציה מסוג `"עיגול לערך הקרוב ביותר`", אפילו לאפשר שימוש בגדלים גדולים יותר של בלוקים בקוונטיזציה (עבורם מחושבים קבועי קוונטיזציה). זה מושג על ידי השארת SW מחוץ לתהליך הקוונטיזציה ושחזור ערכיהם לאחר מכן, תוך צמצום ההשפעות השליליות של ערכים קיצוניים אלה על קוונטיזציה של פרמטרים אחרים. על ידי התמודדות עם האתגרים שמציבים חריגי על-כבד, הגישה המוצעת מאפשרת יישום של שיטות קוונטיזציה פשוטות ויעילות יותר, ומקלה על פריסת מודלים בסביבות עם משאבים מוגבלים.</w:t></w:r></w:p>")
$d.Paragraphs(7).Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t xml:space=`"preserve`">עבודה זו יוצרת טיעון חזק ש-SW אינם רק אנומליות מבודדות אלא רכיבים אינטגרליים הממלאים תפקיד חיוני בעיצוב ההתנהגות והיעילות של LLMs, עם השלכות משמעותיות לדחיסה ולאינפרנס של מודלים. תרומת המאמר אינה טמונה רק בזיהוי SW אלא גם באפיון תפקידם הפונקציונלי בתוך LLMs. המחברים מנתחים כיצד משקלים משפיעים על פלט המודל, ומקשרים אותם ל״התפשטות״ של אקטיבציות חריגות.  </w:t></w:r></w:p>")

# --- Delete the paragraph that starts with "טוב, נכון..." (fully removed in the edit) ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("טוב, נכון שהופיעה")) {
        $p.Range.Delete()
        break
    }
}

# --- Last paragraph: swap the arxiv PDF link for the abs link ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("https://arxiv.org/pdf/2503.01776")) {
        $p.Range.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"Normal`"/></w:pPr><w:r><w:t>https://arxiv.org/abs/2411.07191</w:t></w:r></w:p>")
        break
    }
}

Write-Output ("FinalParaCount=" + $d.Paragraphs.Count)
